$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.339.14'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.403.80'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.62%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +1.46%  '
$ws.Range('E9').Value = '  +8.43%  '
$ws.Range('E10').Value = '  +2.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.73'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.28%  '
$ws.Range('E12').Value = '  +4.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '682.41'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.947.97'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.459.88'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.398.34'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.17%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.121'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.81'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.909'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.88%  '
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.23'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.19%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '103.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.94'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.74'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.75'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.82'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.99'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.17'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '557.47'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('E33').Value = '  +10.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.107'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.45'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.15%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.678.91'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('E38').Value = '  +5.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.83'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0723'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.340'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0426'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.93%  '
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.24%  '
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('E48').Value = '  +5.29%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.53'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.38%  '
$ws.Range('E51').Value = '  +2.81%  '
